# Update nomes integrantes do grupo
#
# The "team members" list (paragraphs styled "Ttulo2") is shifted down by
# one slot and a new member ("Wilbert de Oliveira") is appended, with a
# trailing blank paragraph kept after the list. The "_GoBack" bookmark
# that Word leaves at the last edit position is relocated accordingly.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyInner + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) "*Cliente" paragraph -> single run "<tab/>Daniel Henrique"
#    (the stray tab-only run and the "*Cliente" run collapse into one).
# ---------------------------------------------------------------------
$target = $d.Content.Find
$p = $d.Paragraphs.Item(19)
$full = $p.Range
$body = $d.Range($full.Start, $full.End - 1)
$body.InsertXML((New-PkgXml '<w:body><w:p><w:r><w:rPr><w:i/><w:color w:val="BE994E"/></w:rPr><w:tab/><w:t>Daniel Henrique</w:t></w:r></w:p></w:body>'))

# ---------------------------------------------------------------------
# 2) "Daniel Henrique" -> "Guilherme Gomes"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(20)
$full = $p.Range
$body = $d.Range($full.Start, $full.End - 1)
$body.InsertXML((New-PkgXml '<w:body><w:p><w:r><w:rPr><w:i/><w:color w:val="BE994E"/></w:rPr><w:tab/><w:t>Guilherme Gomes</w:t></w:r></w:p></w:body>'))

# ---------------------------------------------------------------------
# 3) "Guilherme Gomes" -> "Rony Freitas"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(21)
$full = $p.Range
$body = $d.Range($full.Start, $full.End - 1)
$body.InsertXML((New-PkgXml '<w:body><w:p><w:r><w:rPr><w:i/><w:color w:val="BE994E"/></w:rPr><w:tab/><w:t>Rony Freitas</w:t></w:r></w:p></w:body>'))

# ---------------------------------------------------------------------
# 4) "Rony Freitas" -> "Vitor Soares"; paragraph gains a left indent
#    equal to the old tab stop and loses the explicit tab run.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(22)
$p.LeftIndent = 304.8
$full = $p.Range
$body = $d.Range($full.Start, $full.End - 1)
$body.InsertXML((New-PkgXml '<w:body><w:p><w:r><w:rPr><w:i/><w:color w:val="BE994E"/></w:rPr><w:t>Vitor Soares</w:t></w:r></w:p></w:body>'))

# ---------------------------------------------------------------------
# 5) "Vitor Soares" -> "Wilbert de Oliveira"; same indent/tab change as
#    above, plus the "_GoBack" bookmark now lands on this paragraph.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(23)
$p.LeftIndent = 304.8
$full = $p.Range
$bookmarkPoint = $d.Range($full.Start, $full.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)
$full = $p.Range
$body = $d.Range($full.Start, $full.End - 1)
$body.InsertXML((New-PkgXml '<w:body><w:p><w:r><w:rPr><w:i/><w:color w:val="BE994E"/></w:rPr><w:t>Wilbert de Oliveira</w:t></w:r></w:p></w:body>'))

# ---------------------------------------------------------------------
# 6) A new, empty "Ttulo2" paragraph follows the (now renamed) last
#    member of the list.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(23)
$endOfPara = $d.Range($p.Range.End, $p.Range.End)
$endOfPara.InsertXML((New-PkgXml '<w:body><w:p><w:pPr><w:pStyle w:val="Ttulo2"/><w:tabs><w:tab w:val="left" w:pos="6096"/></w:tabs><w:spacing w:before="120"/><w:rPr><w:i/><w:color w:val="BE994E"/></w:rPr></w:pPr></w:p></w:body>'))

# Note: Word keeps a single "_GoBack" bookmark; re-adding it under (5)
# above already relocated it away from the old "GitHub" paragraph, so
# no separate removal step is needed here.
